# Applies the scheduled runner updates to the Ultima_Profits leve-crafting profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 814.8200000000001
$ws.Range("I15").Value = 814.8200000000001
$ws.Range("K15").Value = 2444.46
$ws.Range("M15").Value = -2275.46

# Row 19
$ws.Range("H19").Value = 1103.6207
$ws.Range("I19").Value = 1082.8667
$ws.Range("J19").Value = 1125.8572
$ws.Range("K19").Value = 1082.8667
$ws.Range("L19").Value = 1125.8572
$ws.Range("M19").Value = -907.8667
$ws.Range("N19").Value = -1475.8572

# Row 53
$ws.Range("H53").Value = 1309
$ws.Range("I53").Value = 3700
$ws.Range("J53").Value = 113.5
$ws.Range("K53").Value = 3700
$ws.Range("L53").Value = 113.5
$ws.Range("M53").Value = -3063
$ws.Range("N53").Value = -1387.5

# Row 55
$ws.Range("H55").Value = 162.47058
$ws.Range("I55").Value = 156.2
$ws.Range("J55").Value = 165.08333
$ws.Range("K55").Value = 156.2
$ws.Range("L55").Value = 165.08333
$ws.Range("M55").Value = 57.80000000000001
$ws.Range("N55").Value = -593.0833299999999

# Row 138
$ws.Range("H138").Value = 1700.5714
$ws.Range("I138").Value = 1577.4333
$ws.Range("J138").Value = 2439.4
$ws.Range("K138").Value = 4732.2999
$ws.Range("L138").Value = 7318.200000000001
$ws.Range("M138").Value = 407.7001
$ws.Range("N138").Value = -17598.2

$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 30100
$ws.Range("J8").Value = 30100
$ws.Range("L8").Value = 30100
$ws.Range("N8").Value = -30388

# Row 32
$ws.Range("H32").Value = 12488.462
$ws.Range("I32").Value = 12549.048
$ws.Range("K32").Value = 12549.048
$ws.Range("M32").Value = -12262.048

# Row 43
$ws.Range("H43").Value = 9761.666999999999
$ws.Range("J43").Value = 9501.625
$ws.Range("L43").Value = 9501.625
$ws.Range("N43").Value = -10127.625

# Row 46
$ws.Range("H46").Value = 30000
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1179.25
$ws.Range("I94").Value = 1020.6667
$ws.Range("J94").Value = 1655
$ws.Range("K94").Value = 1020.6667
$ws.Range("L94").Value = 1655
$ws.Range("M94").Value = -569.6667
$ws.Range("N94").Value = -2557

# Row 134
$ws.Range("H134").Value = 3160.5095
$ws.Range("I134").Value = 2316.275
$ws.Range("J134").Value = 5758.154
$ws.Range("K134").Value = 6948.825000000001
$ws.Range("L134").Value = 17274.462
$ws.Range("M134").Value = -4413.825000000001
$ws.Range("N134").Value = -22344.462

$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 351670.66
$ws.Range("I6").Value = 351670.66
$ws.Range("K6").Value = 351670.66
$ws.Range("M6").Value = -351557.66

# Row 31
$ws.Range("H31").Value = 6064376.5
$ws.Range("I31").Value = 4110.3413
$ws.Range("J31").Value = 23812298
$ws.Range("K31").Value = 4110.3413
$ws.Range("L31").Value = 23812298
$ws.Range("M31").Value = -3815.3413
$ws.Range("N31").Value = -23812888

# Row 34
$ws.Range("H34").Value = 6064376.5
$ws.Range("I34").Value = 4110.3413
$ws.Range("J34").Value = 23812298
$ws.Range("K34").Value = 4110.3413
$ws.Range("L34").Value = 23812298
$ws.Range("M34").Value = -3908.3413
$ws.Range("N34").Value = -23812702

# Row 134
$ws.Range("H134").Value = 1014.26154
$ws.Range("I134").Value = 922.45
$ws.Range("J134").Value = 2116
$ws.Range("K134").Value = 2767.35
$ws.Range("L134").Value = 6348
$ws.Range("M134").Value = -232.3500000000004
$ws.Range("N134").Value = -11418

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 57.25
$ws.Range("I12").Value = 51.2
$ws.Range("J12").Value = 61.57143
$ws.Range("K12").Value = 153.6
$ws.Range("L12").Value = 184.71429
$ws.Range("M12").Value = 19.39999999999998
$ws.Range("N12").Value = -530.71429

# Row 22
$ws.Range("H22").Value = 1074.9166
$ws.Range("I22").Value = 483.33334
$ws.Range("J22").Value = 1666.5
$ws.Range("K22").Value = 1450.00002
$ws.Range("L22").Value = 4999.5
$ws.Range("M22").Value = -1281.00002
$ws.Range("N22").Value = -5337.5

# Row 23
$ws.Range("H23").Value = 130.70589
$ws.Range("I23").Value = 24.4
$ws.Range("J23").Value = 175
$ws.Range("K23").Value = 73.19999999999999
$ws.Range("L23").Value = 525
$ws.Range("M23").Value = 161.8
$ws.Range("N23").Value = -995

# Row 27
$ws.Range("H27").Value = 1074.9166
$ws.Range("I27").Value = 483.33334
$ws.Range("J27").Value = 1666.5
$ws.Range("K27").Value = 1450.00002
$ws.Range("L27").Value = 4999.5
$ws.Range("M27").Value = -1348.00002
$ws.Range("N27").Value = -5203.5

# Row 131
$ws.Range("H131").Value = 848.2
$ws.Range("I131").Value = 466.66666
$ws.Range("J131").Value = 860
$ws.Range("K131").Value = 1399.99998
$ws.Range("L131").Value = 2580
$ws.Range("M131").Value = 3640.00002
$ws.Range("N131").Value = -12660

# Row 132
$ws.Range("H132").Value = 2345
$ws.Range("I132").Value = 772.5
$ws.Range("K132").Value = 6952.5
$ws.Range("M132").Value = -4422.5

$ws = $wb.Worksheets.Item("GSM")
# Row 119
$ws.Range("H119").Value = 37561
$ws.Range("J119").Value = 37561
$ws.Range("L119").Value = 37561
$ws.Range("N119").Value = -47237

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1640.6
$ws.Range("I16").Value = 1835.5385
$ws.Range("K16").Value = 1835.5385
$ws.Range("M16").Value = -1665.5385

# Row 40
$ws.Range("H40").Value = 12312.714
$ws.Range("I40").Value = 14547.5
$ws.Range("J40").Value = 9333
$ws.Range("K40").Value = 14547.5
$ws.Range("L40").Value = 9333
$ws.Range("M40").Value = -14411.5
$ws.Range("N40").Value = -9605

# Row 46
$ws.Range("H46").Value = 684.45
$ws.Range("I46").Value = 676.125
$ws.Range("J46").Value = 690
$ws.Range("K46").Value = 676.125
$ws.Range("L46").Value = 690
$ws.Range("M46").Value = -488.125
$ws.Range("N46").Value = -1066

# Row 122
$ws.Range("H122").Value = 4916.9165
$ws.Range("I122").Value = 4504.037
$ws.Range("K122").Value = 13512.111
$ws.Range("M122").Value = -11062.111

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 10280.667
$ws.Range("I62").Value = 6134.9
$ws.Range("K62").Value = 6134.9
$ws.Range("M62").Value = -5510.9

# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 65
$ws.Range("H65").Value = 10280.667
$ws.Range("I65").Value = 6134.9
$ws.Range("K65").Value = 30674.5
$ws.Range("M65").Value = -27554.5

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Row 96
$ws.Range("H96").Value = 3706.6365
$ws.Range("I96").Value = 3136.8333
$ws.Range("J96").Value = 4390.4
$ws.Range("K96").Value = 3136.8333
$ws.Range("L96").Value = 4390.4
$ws.Range("M96").Value = -1763.8333
$ws.Range("N96").Value = -7136.4

# Row 132
$ws.Range("H132").Value = 1323.2941
$ws.Range("I132").Value = 952.3051
$ws.Range("K132").Value = 2856.9153
$ws.Range("M132").Value = -326.9153000000001
